# KCOR_Bias_Correction_Table.xlsx edit
#
# - bold the header row labels (A1 empty cell + B1:F1 "1x".."5x"),
#   with B1:F1 additionally right-aligned
# - remove the long explanatory sentence that used to sit in I1
# - re-enter that sentence, split in two, in A11/A12 instead
# - leave the "source" (I3) / URL (I4) cells as-is (their shared-string
#   index shifts automatically once the old string is removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bold the (empty) A1 cell and the B1:F1 header labels; right-align B1:F1
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1:F1").Font.Bold = $true
$ws.Range("B1:F1").HorizontalAlignment = -4152  # xlRight

# Drop the old long note that lived in I1
$ws.Range("I1").ClearContents()

# Re-add the note, split across two rows, near the bottom of the table
$ws.Range("A11").Value = "this is the annual bias correction for cohorts with enhanced frailty"
$ws.Range("A12").Value = " (measured by ACM ratios at that age)"

# Restore last-used selection
$ws.Range("G17").Select() | Out-Null
